# Update "Forecast Comparison" sheet with a new Week_Start_Date column,
# corrected MyForecast values, shorter week labels, and a proper boolean
# is_holiday_week column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the old column B (ASIN), so everything from
# ASIN onward shifts one column to the right.
$ws.Columns.Item(2).Insert()

# New header + per-row Week_Start_Date values.
$ws.Range("B1").Value = "Week_Start_Date"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]
    $cell.Style = "Normal"
}

# Shorten the week labels from "W01".."W16" to "W1".."W16".
for ($i = 1; $i -le 16; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = "W$i"
}

# The old "is_holiday_week" column (now column J after the insert) holds
# plain numbers; convert it to a proper boolean FALSE value.
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}

$wb.Save()
